# Fix purchase/sales order resource/table/fields.
# Rebuilds Sheet1 with the new header/column layout and new product rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- wipe the old table (values + formatting) -----------------------------
$ws.Cells.ClearContents()
$ws.Cells.ClearFormats()

# --- header row -------------------------------------------------------------
$ws.Range("A1").Value = "名称"
$ws.Range("B1").Value = "分类"
$ws.Range("C1").Value = "供应商"
$ws.Range("D1").Value = "品牌"
$ws.Range("E1").Value = "价格"
$ws.Range("F1").Value = "原价"
$ws.Range("G1").Value = "采购价"
$ws.Range("H1").Value = "库存"
$ws.Range("I1").Value = "状态"
$ws.Range("J1").Value = "详情"

# highlighted header cells (red font)
$ws.Range("A1").Font.Color = 255
$ws.Range("B1").Font.Color = 255
$ws.Range("E1").Font.Color = 255

# --- data rows ---------------------------------------------------------------
$ws.Range("A2").Value = "雀巢怡养中老年奶粉"
$ws.Range("B2").Value = "营养"
$ws.Range("C2").Value = "新新健康"
$ws.Range("D2").Value = "雀巢"
$ws.Range("E2").Value = 178
$ws.Range("F2").Value = 239
$ws.Range("G2").Value = 88
$ws.Range("H2").Value = 20
$ws.Range("I2").Value = "上架"

$ws.Range("A3").Value = "康恩贝维生素vc咀嚼片"
$ws.Range("B3").Value = "营养"
$ws.Range("C3").Value = "新新健康"
$ws.Range("D3").Value = "康恩贝"
$ws.Range("E3").Value = 19.9
$ws.Range("F3").Value = 32.8
$ws.Range("G3").Value = 6.8
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = "上架"

$ws.Range("A4").Value = "本博颈椎按摩器按摩枕头"
$ws.Range("B4").Value = "保健"
$ws.Range("C4").Value = "新新健康"
$ws.Range("D4").Value = "本博"
$ws.Range("E4").Value = 0.02
$ws.Range("F4").Value = 338
$ws.Range("G4").Value = 0.01
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = "上架"

$ws.Range("A5").Value = "中老年澳洲复合燕麦片700g"
$ws.Range("B5").Value = "营养"
$ws.Range("C5").Value = "新新健康"
$ws.Range("D5").Value = "西麦"
$ws.Range("E5").Value = 0.01
$ws.Range("F5").Value = 59
$ws.Range("G5").Value = 0.01
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = "推荐"

# "详情" column cells are wrapped but left blank (template placeholder)
$ws.Range("J2").WrapText = $true
$ws.Range("J3").WrapText = $true
$ws.Range("J5").WrapText = $true

# row 3 no longer needs the taller custom height from the old layout
$ws.Rows.Item(3).AutoFit()

# --- column sizing -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 29.833333333333336
$ws.Columns.Item(10).ColumnWidth = 14.666666666666666

# --- selection (matches the saved cursor position in the workbook) ---------
$ws.Range("B8").Select()
